$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "changed batch size to 16": add the MSE (column D) figures for the
# LR1Batch and LR3Batch rows (batch size 16 results) in both the
# "Comparison of batch size" table and the "Augmented vs. Unaugmented"
# table, plus the Unaligned 200-epoch MSE value used in two tables.
$ws.Range("D15").Value = 1032.94294633
$ws.Range("D16").Value = 1095.0160229200001
$ws.Range("D26").Value = 10669.743807999999
$ws.Range("D32").Value = 1095.0160229200001
$ws.Range("D33").Value = 1032.94294633
$ws.Range("D37").Value = 10669.743807999999

# Reflect the updated view state (scroll position / zoom / active cell).
$excel.ActiveWindow.Zoom = 85
[void]$ws.Range("A44").Select()
